# Applies the cryptos price/volume refresh described by the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some refreshed "Price" values are purely numeric-looking strings
# (e.g. "200.79"). The sheet stores Price/Volume as plain text, so we
# pre-format those specific cells as Text before assigning, then reset
# the style back to Normal (keeps the value as text without leaving a
# stray direct-formatted style behind).
$numericPriceCells = @("D5", "D6", "D15", "D19", "D21", "D22", "D25", "D26", "D27", "D28", "D30", "D31", "D33", "D34", "D36", "D37", "D39", "D40", "D41", "D44", "D45", "D48", "D51")
foreach ($cellRef in $numericPriceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "76.503.17"
$ws.Range("E2").Value = "  +0.75%  "

$ws.Range("D3").Value = "3.037.26"
$ws.Range("E3").Value = "  +4.19%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "200.79"
$ws.Range("E5").Value = "  +0.92%  "

$ws.Range("D6").Value = "631.53"
$ws.Range("E6").Value = "  +5.52%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E8").Value = "  +0.15%  "

$ws.Range("E9").Value = "  +1.80%  "

$ws.Range("D10").Value = "3.035.82"
$ws.Range("E10").Value = "  +4.17%  "

$ws.Range("E11").Value = "  -0.76%  "

$ws.Range("E12").Value = "  -0.42%  "

$ws.Range("D14").Value = "3.597.29"
$ws.Range("E14").Value = "  +4.25%  "

$ws.Range("D15").Value = "29.36"
$ws.Range("E15").Value = "  +6.78%  "

$ws.Range("D16").Value = "76.408.29"
$ws.Range("E16").Value = "  +0.76%  "

$ws.Range("E17").Value = "  -1.07%  "

$ws.Range("D18").Value = "3.042.26"
$ws.Range("E18").Value = "  +4.39%  "

$ws.Range("D19").Value = "13.50"
$ws.Range("E19").Value = "  +5.40%  "

$ws.Range("E20").Value = "  +1.62%  "

$ws.Range("D21").Value = "375.37"
$ws.Range("E21").Value = "  -0.97%  "

$ws.Range("D22").Value = "4.35"
$ws.Range("E22").Value = "  +1.97%  "

$ws.Range("E23").Value = "  -1.52%  "

$ws.Range("D24").Value = "3.206.35"
$ws.Range("E24").Value = "  +4.61%  "

$ws.Range("D25").Value = "73.05"
$ws.Range("E25").Value = "  +2.52%  "

$ws.Range("D26").Value = "4.39"
$ws.Range("E26").Value = "  +4.37%  "

$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  +0.01%  "

$ws.Range("D28").Value = "9.83"
$ws.Range("E28").Value = "  +1.71%  "

$ws.Range("E29").Value = "  +0.25%  "

$ws.Range("D30").Value = "0.995"
$ws.Range("E30").Value = "  -0.38%  "

$ws.Range("D31").Value = "8.34"
$ws.Range("E31").Value = "  +7.84%  "

$ws.Range("E32").Value = "  -1.36%  "

$ws.Range("D33").Value = "512.91"
$ws.Range("E33").Value = "  +1.64%  "

$ws.Range("D34").Value = "1.95"
$ws.Range("E34").Value = "  +6.69%  "

$ws.Range("E35").Value = "  -0.02%  "

$ws.Range("D36").Value = "20.84"
$ws.Range("E36").Value = "  +2.90%  "

$ws.Range("D37").Value = "164.26"
$ws.Range("E37").Value = "  -0.31%  "

$ws.Range("E38").Value = "  +1.71%  "

$ws.Range("D39").Value = "0.382"
$ws.Range("E39").Value = "  +10.73%  "

$ws.Range("D40").Value = "192.09"
$ws.Range("E40").Value = "  +6.39%  "

$ws.Range("D41").Value = "0.105"
$ws.Range("E41").Value = "  -2.07%  "

$ws.Range("E42").Value = "  -1.65%  "

$ws.Range("E43").Value = "  +0.19%  "

$ws.Range("D44").Value = "5.06"
$ws.Range("E44").Value = "  +1.18%  "

$ws.Range("D45").Value = "42.75"
$ws.Range("E45").Value = "  +6.16%  "

$ws.Range("E46").Value = "  +4.32%  "

$ws.Range("E47").Value = "  -0.30%  "

$ws.Range("D48").Value = "0.607"
$ws.Range("E48").Value = "  +6.20%  "

$ws.Range("E49").Value = "  +5.84%  "

$ws.Range("E50").Value = "  +1.33%  "

$ws.Range("D51").Value = "3.89"
$ws.Range("E51").Value = "  +4.16%  "

foreach ($cellRef in $numericPriceCells) {
    $ws.Range($cellRef).Style = "Normal"
}
